# RSTK-9586-SYDATA-Work Order Close.xlsx
#
# The second data row (row 3) was a duplicate of row 2 except that its
# "Background Processing" flag (column G) was TRUE while row 2's was
# FALSE. The edit removes that duplicate row and flips row 2's flag to
# TRUE, then leaves the whole (now-last) data row selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Flip "Background Processing" (G2) from FALSE to TRUE
$ws.Range("G2").Value = $true

# Remove the now-redundant row 3 entirely (rows below shift up)
$ws.Rows("3:3").Delete() | Out-Null

# Leave the full row 2 selected, matching the post-edit view state
$ws.Rows("2:2").Select() | Out-Null
